{"js": "// Office.js (Word JavaScript API) script.\n//\n// Applies two edits to the document:\n//\n//  1) The paragraph holding the runs \"{\" \"Image1\" \"}\" (a placeholder\n//     token, Calibri 11pt) is split into two paragraphs:\n//       - the same paragraph, now styled Segoe UI 11.5pt/black and\n//         holding two runs \"ACE_\" + \"Image1\" (text \"ACE_Image1\")\n//       - a brand-new paragraph right after it that keeps the original\n//         Calibri 11pt styling and the original token text \"{Image1}\"\n//         (now a single run)\n//\n//  2) The paragraph \"i hope you are alright\" is split into two runs,\n//     \"i\" and \" hope you are alright\", with proofing-error bookmarks\n//     (spellStart/spellEnd) wrapped around the first run \"i\" (as Word's\n//     spell checker would mark a lowercase standalone \"i\").\n\nconst OOXML_WRAPPER =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>{BODY}</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nfunction wrapOoxml(bodyInner) {\n  return OOXML_WRAPPER.replace(\"{BODY}\", bodyInner);\n}\n\n// Pull the opening \"<w:p ...>\" tag (with its paraId/rsid attributes, if\n// any) straight out of a paragraph's own OOXML so the replacement keeps\n// the same paragraph identity Word assigned it.\nasync function openingParaTag(paragraph) {\n  const result = paragraph.getOoxml();\n  await paragraph.context.sync();\n  const match = /<w:p(\\s[^>]*)?>/.exec(result.value);\n  return match ? \"<w:p\" + (match[1] || \"\") + \">\" : \"<w:p>\";\n}\n\nasync function findParagraphByText(body, text) {\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === text) {\n      return paragraphs.items[i];\n    }\n  }\n  return null;\n}\n\nconst body = context.document.body;\n\n// --- Edit 1: \"{Image1}\" placeholder paragraph -> \"ACE_Image1\" (Segoe UI)\n//     followed by a new \"{Image1}\" paragraph (original Calibri style).\nconst tokenParagraph = await findParagraphByText(body, \"{Image1}\");\nif (tokenParagraph) {\n  const openingTag = await openingParaTag(tokenParagraph);\n  const tokenRange = tokenParagraph.getRange();\n  const replacementOoxml =\n    openingTag +\n    '<w:pPr><w:spacing w:after=\"200\" w:line=\"276\" w:lineRule=\"auto\"/>' +\n    '<w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n    '<w:color w:val=\"000000\"/><w:sz w:val=\"23\"/><w:szCs w:val=\"23\"/></w:rPr>' +\n    \"</w:pPr>\" +\n    '<w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n    '<w:color w:val=\"000000\"/><w:sz w:val=\"23\"/><w:szCs w:val=\"23\"/></w:rPr>' +\n    \"<w:t>ACE_</w:t></w:r>\" +\n    '<w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n    '<w:color w:val=\"000000\"/><w:sz w:val=\"23\"/><w:szCs w:val=\"23\"/></w:rPr>' +\n    \"<w:t>Image1</w:t></w:r>\" +\n    \"</w:p>\" +\n    \"<w:p>\" +\n    '<w:pPr><w:spacing w:after=\"200\" w:line=\"276\" w:lineRule=\"auto\"/>' +\n    '<w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"22\"/>' +\n    '<w:lang w:val=\"en\"/></w:rPr>' +\n    \"</w:pPr>\" +\n    '<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"22\"/>' +\n    '<w:lang w:val=\"en\"/></w:rPr>' +\n    \"<w:t>{Image1}</w:t></w:r>\" +\n    \"</w:p>\";\n  tokenRange.insertOoxml(wrapOoxml(replacementOoxml), \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 2: \"i hope you are alright\" -> split into \"i\" / \" hope you\n//     are alright\" runs, flanked by spellStart/spellEnd proofing marks.\n// Re-query the paragraph collection since edit 1 changed paragraph\n// count/indices and any proxies captured beforehand would be stale.\nconst greetingParagraph = await findParagraphByText(\n  body,\n  \"i hope you are alright\"\n);\nif (greetingParagraph) {\n  const openingTag = await openingParaTag(greetingParagraph);\n  const greetingRange = greetingParagraph.getRange();\n  const replacementOoxml =\n    openingTag +\n    '<w:pPr><w:spacing w:after=\"200\" w:line=\"276\" w:lineRule=\"auto\"/>' +\n    '<w:rPr><w:lang w:val=\"en\"/></w:rPr>' +\n    \"</w:pPr>\" +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"22\"/>' +\n    '<w:lang w:val=\"en\"/></w:rPr>' +\n    \"<w:t>i</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"22\"/>' +\n    '<w:lang w:val=\"en\"/></w:rPr>' +\n    '<w:t xml:space=\"preserve\"> hope you are alright</w:t></w:r>' +\n    \"</w:p>\";\n  greetingRange.insertOoxml(wrapOoxml(replacementOoxml), \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n#\n# Applies two edits to the document:\n#\n#  1) The paragraph holding the runs \"{\" \"Image1\" \"}\" (a placeholder\n#     token, Calibri 11pt) is split into two paragraphs:\n#       - the same paragraph, now styled Segoe UI 11.5pt/black and\n#         holding two runs \"ACE_\" + \"Image1\" (text \"ACE_Image1\")\n#       - a brand-new paragraph right after it that keeps the original\n#         Calibri 11pt styling and the original token text \"{Image1}\"\n#         (now a single run)\n#\n#  2) The paragraph \"i hope you are alright\" is split into two runs,\n#     \"i\" and \" hope you are alright\", with proofing-error bookmarks\n#     (spellStart/spellEnd) wrapped around the first run \"i\" (as Word's\n#     spell checker would mark a lowercase standalone \"i\").\n\n$d = $word.ActiveDocument\n\nfunction Wrap-Ooxml([string]$bodyInner) {\n    return '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\nfunction Get-OpeningParaTag($paragraph) {\n    $xml = $paragraph.Range.WordOpenXML\n    if ($xml -match '<w:p(\\s[^>]*)?>') {\n        return $matches[0]\n    }\n    return '<w:p>'\n}\n\nfunction Find-ParagraphByText($doc, [string]$text) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text -eq ($text + \"`r\")) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# --- Edit 1: \"{Image1}\" placeholder paragraph -> \"ACE_Image1\" (Segoe UI)\n#     followed by a new \"{Image1}\" paragraph (original Calibri style).\n$tokenParagraph = Find-ParagraphByText $d \"{Image1}\"\nif ($tokenParagraph -ne $null) {\n    $openingTag = Get-OpeningParaTag $tokenParagraph\n    $replacement = $openingTag +\n        '<w:pPr><w:spacing w:after=\"200\" w:line=\"276\" w:lineRule=\"auto\"/>' +\n        '<w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n        '<w:color w:val=\"000000\"/><w:sz w:val=\"23\"/><w:szCs w:val=\"23\"/></w:rPr>' +\n        '</w:pPr>' +\n        '<w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n        '<w:color w:val=\"000000\"/><w:sz w:val=\"23\"/><w:szCs w:val=\"23\"/></w:rPr>' +\n        '<w:t>ACE_</w:t></w:r>' +\n        '<w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n        '<w:color w:val=\"000000\"/><w:sz w:val=\"23\"/><w:szCs w:val=\"23\"/></w:rPr>' +\n        '<w:t>Image1</w:t></w:r>' +\n        '</w:p>' +\n        '<w:p>' +\n        '<w:pPr><w:spacing w:after=\"200\" w:line=\"276\" w:lineRule=\"auto\"/>' +\n        '<w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"22\"/>' +\n        '<w:lang w:val=\"en\"/></w:rPr>' +\n        '</w:pPr>' +\n        '<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"22\"/>' +\n        '<w:lang w:val=\"en\"/></w:rPr>' +\n        '<w:t>{Image1}</w:t></w:r>' +\n        '</w:p>'\n    $tokenParagraph.Range.InsertXML((Wrap-Ooxml $replacement))\n}\n\n# --- Edit 2: \"i hope you are alright\" -> split into \"i\" / \" hope you\n#     are alright\" runs, flanked by spellStart/spellEnd proofing marks.\n# Re-scan the paragraph collection: edit 1 above changed the paragraph\n# count, so paragraph objects must be looked up again by text.\n$greetingParagraph = Find-ParagraphByText $d \"i hope you are alright\"\nif ($greetingParagraph -ne $null) {\n    $openingTag = Get-OpeningParaTag $greetingParagraph\n    $replacement = $openingTag +\n        '<w:pPr><w:spacing w:after=\"200\" w:line=\"276\" w:lineRule=\"auto\"/>' +\n        '<w:rPr><w:lang w:val=\"en\"/></w:rPr>' +\n        '</w:pPr>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"22\"/>' +\n        '<w:lang w:val=\"en\"/></w:rPr>' +\n        '<w:t>i</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"22\"/>' +\n        '<w:lang w:val=\"en\"/></w:rPr>' +\n        '<w:t xml:space=\"preserve\"> hope you are alright</w:t></w:r>' +\n        '</w:p>'\n    $greetingParagraph.Range.InsertXML((Wrap-Ooxml $replacement))\n}\n"}
